# Daily scrape update - 2025-08-11 03:42:05 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# (ColumnWidth setter rounds to 1/6-pixel steps with a ~0.83 char offset vs.
#  the stored OOXML <col width> value, so subtract 0.83 from the desired
#  stored width to land exactly on the target.)
$ws.Columns.Item(3).ColumnWidth = 58.17   # stored width 23 -> 59
$ws.Columns.Item(4).ColumnWidth = 22.17   # stored width 70 -> 23
$ws.Columns.Item(6).ColumnWidth = 16.17   # stored width 15 -> 17
$ws.Columns.Item(7).ColumnWidth = 15.17   # stored width 15 -> 16
$ws.Columns.Item(8).ColumnWidth = 24.17   # stored width 16 -> 25

# --- Row 2: new opportunity record ---
$ws.Range("A2").Value = "'1325417"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1325417"
$ws.Range("C2").Value = "Junior Software Engineer – AI & Internal Tools (EU ONLY)"
$ws.Range("D2").Value = "Brussels, Belgium"
$ws.Range("F2").Value = "69 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "Eureka Resource Mining"

# --- Row 3: new opportunity record ---
$ws.Range("A3").Value = "'1315734"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1315734"
$ws.Range("C3").Value = "Accelerate Romania - WordPress Web development"
$ws.Range("D3").Value = "Târgu Mureș, Romania"
$ws.Range("F3").Value = "135 applicants"
$ws.Range("H3").Value = "Streamline Media"
